$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it at the end -
# renaming sheets and nudging the selection on another sheet must not change
# which tab is active in the saved workbook.
$originalActiveSheetName = $wb.ActiveSheet.Name

# Rename "wt" -> "wt_log2_expression" and "dcin5" -> "dcin5_log2_expression"
# (inputs now carry the "_log2_expression" suffix).
$wb.Worksheets.Item("wt").Name = "wt_log2_expression"
$wb.Worksheets.Item("dcin5").Name = "dcin5_log2_expression"

# Move the selection on the (renamed) dcin5 sheet from O9 to F43. Excel only
# lets you change the selection of the sheet that is currently active, so we
# briefly activate it, move the selection, then reactivate the sheet that was
# active originally.
$dcin5 = $wb.Worksheets.Item("dcin5_log2_expression")
$dcin5.Activate()
$dcin5.Range("F43").Select()

$wb.Worksheets.Item($originalActiveSheetName).Activate()
